$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"23.9300000000003"
$ws.Range("L2").Value = [double]"62.28997642295737"
$ws.Range("M2").Value = "[53.046088207082775, 71.53386463883196]"
$ws.Range("N2").Value = [double]"0"
$ws.Range("O2").Value = [double]"0"
$ws.Range("P2").Value = [double]"1.490605523324887"
$ws.Range("Q2").Value = "[1.3270791789938858, 1.654131867655888]"
$ws.Range("T2").Value = [double]"56.41309856561734"
$ws.Range("U2").Value = "[50.45677299271063, 62.369424138524046]"
$ws.Range("X2").Value = [double]"18.25291291291314"
$ws.Range("Y2").Value = [double]"17.63011011011033"
$ws.Range("Z2").Value = [double]"18.87571571571596"
$ws.Range("F3").Value = [double]"23.9300000000003"
$ws.Range("H3").Value = [double]"4.101818884549857e-11"
$ws.Range("I3").Value = [double]"4.101818884549857e-11"
$ws.Range("L3").Value = [double]"55.04658486227049"
$ws.Range("M3").Value = "[41.34909143231873, 68.74407829222226]"
$ws.Range("N3").Value = [double]"2.515403441094577e-10"
$ws.Range("O3").Value = [double]"2.515403441094577e-10"
$ws.Range("P3").Value = [double]"1.691868716347656"
$ws.Range("Q3").Value = "[1.4151318259413488, 1.9686056067539637]"
$ws.Range("R3").Value = [double]"4.440892098500626e-16"
$ws.Range("S3").Value = [double]"4.440892098500626e-16"
$ws.Range("T3").Value = [double]"54.86484130182046"
$ws.Range("U3").Value = "[46.38895154876326, 63.34073105487765]"
$ws.Range("X3").Value = [double]"17.48638638638661"
$ws.Range("Y3").Value = [double]"16.43241241241262"
$ws.Range("Z3").Value = [double]"18.54036036036059"
$ws.Range("F4").Value = [double]"23.9300000000003"
$ws.Range("H4").Value = [double]"2.645705876602733e-11"
$ws.Range("I4").Value = [double]"2.645705876602733e-11"
$ws.Range("L4").Value = [double]"64.3872930839208"
$ws.Range("M4").Value = "[49.062707155127754, 79.71187901271384]"
$ws.Range("N4").Value = [double]"7.391265377520995e-11"
$ws.Range("O4").Value = [double]"7.391265377520995e-11"
$ws.Range("P4").Value = [double]"1.616395018964117"
$ws.Range("Q4").Value = "[1.352237078121731, 1.8805529598065034]"
$ws.Range("R4").Value = [double]"4.440892098500626e-16"
$ws.Range("S4").Value = [double]"4.440892098500626e-16"
$ws.Range("T4").Value = [double]"58.35319329952882"
$ws.Range("U4").Value = "[48.58413520391218, 68.12225139514545]"
$ws.Range("V4").Value = [double]"1.110223024625157e-15"
$ws.Range("W4").Value = [double]"1.110223024625157e-15"
$ws.Range("X4").Value = [double]"17.77383383383406"
$ws.Range("Y4").Value = [double]"16.76776776776798"
$ws.Range("Z4").Value = [double]"18.77989989990014"
$ws.Range("F5").Value = [double]"23.9300000000003"
$ws.Range("H5").Value = [double]"2.886579864025407e-15"
$ws.Range("I5").Value = [double]"2.886579864025407e-15"
$ws.Range("L5").Value = [double]"61.81810231762874"
$ws.Range("M5").Value = "[50.30479420865936, 73.33141042659813]"
$ws.Range("N5").Value = [double]"4.218847493575595e-14"
$ws.Range("O5").Value = [double]"4.218847493575595e-14"
$ws.Range("P5").Value = [double]"1.717026615475502"
$ws.Range("Q5").Value = "[1.515763422452732, 1.9182898084982725]"
$ws.Range("T5").Value = [double]"52.37050602033886"
$ws.Range("U5").Value = "[45.29958757048685, 59.44142447019087]"
$ws.Range("X5").Value = [double]"17.39057057057079"
$ws.Range("Y5").Value = [double]"16.62404404404425"
$ws.Range("Z5").Value = [double]"18.15709709709733"
$ws.Range("F6").Value = [double]"23.9300000000003"
$ws.Range("H6").Value = [double]"3.906874823655926e-13"
$ws.Range("I6").Value = [double]"3.906874823655926e-13"
$ws.Range("L6").Value = [double]"58.17853147594043"
$ws.Range("M6").Value = "[45.77221125891842, 70.58485169296245]"
$ws.Range("N6").Value = [double]"3.015365734881925e-12"
$ws.Range("O6").Value = [double]"3.015365734881925e-12"
$ws.Range("P6").Value = [double]"1.729605565039425"
$ws.Range("Q6").Value = "[1.478026573760963, 1.9811845563178876]"
$ws.Range("T6").Value = [double]"58.16325619940233"
$ws.Range("U6").Value = "[50.4608739282069, 65.86563847059776]"
$ws.Range("X6").Value = [double]"17.34266266266288"
$ws.Range("Y6").Value = [double]"16.38450450450471"
$ws.Range("Z6").Value = [double]"18.30082082082105"
$ws.Range("F7").Value = [double]"23.9300000000003"
$ws.Range("H7").Value = [double]"6.372680161348399e-14"
$ws.Range("I7").Value = [double]"6.372680161348399e-14"
$ws.Range("L7").Value = [double]"63.6604548268126"
$ws.Range("M7").Value = "[50.32270689658493, 76.99820275704027]"
$ws.Range("N7").Value = [double]"1.763034163104749e-12"
$ws.Range("O7").Value = [double]"1.763034163104749e-12"
$ws.Range("P7").Value = [double]"1.767342413731195"
$ws.Range("Q7").Value = "[1.54092132158058, 1.9937635058818106]"
$ws.Range("R7").Value = [double]"0"
$ws.Range("S7").Value = [double]"0"
$ws.Range("T7").Value = [double]"59.04687988879965"
$ws.Range("U7").Value = "[51.06482321032534, 67.02893656727396]"
$ws.Range("V7").Value = [double]"0"
$ws.Range("W7").Value = [double]"0"
$ws.Range("X7").Value = [double]"17.19893893893915"
$ws.Range("Y7").Value = [double]"16.3365965965968"
$ws.Range("Z7").Value = [double]"18.06128128128151"
$ws.Range("F8").Value = [double]"23.9300000000003"
$ws.Range("H8").Value = [double]"1.046607245314135e-12"
$ws.Range("I8").Value = [double]"1.046607245314135e-12"
$ws.Range("L8").Value = [double]"62.41984283804997"
$ws.Range("M8").Value = "[48.40546592178147, 76.43421975431848]"
$ws.Range("N8").Value = [double]"1.393418713746541e-11"
$ws.Range("O8").Value = [double]"1.393418713746541e-11"
$ws.Range("P8").Value = [double]"1.729605565039425"
$ws.Range("Q8").Value = "[1.478026573760963, 1.9811845563178876]"
$ws.Range("T8").Value = [double]"54.97726278068576"
$ws.Range("U8").Value = "[46.45784751811383, 63.49667804325768]"
$ws.Range("V8").Value = [double]"0"
$ws.Range("W8").Value = [double]"0"
$ws.Range("X8").Value = [double]"17.34266266266288"
$ws.Range("Y8").Value = [double]"16.38450450450471"
$ws.Range("Z8").Value = [double]"18.30082082082105"
$ws.Range("F9").Value = [double]"23.9300000000003"
$ws.Range("H9").Value = [double]"2.240474472614551e-11"
$ws.Range("I9").Value = [double]"2.240474472614551e-11"
$ws.Range("L9").Value = [double]"57.15405151059129"
$ws.Range("M9").Value = "[42.74579785797964, 71.56230516320295]"
$ws.Range("N9").Value = [double]"3.571001272462127e-10"
$ws.Range("O9").Value = [double]"3.571001272462127e-10"
$ws.Range("P9").Value = [double]"1.767342413731195"
$ws.Range("Q9").Value = "[1.478026573760964, 2.0566582537014266]"
$ws.Range("R9").Value = [double]"6.661338147750939e-16"
$ws.Range("S9").Value = [double]"6.661338147750939e-16"
$ws.Range("T9").Value = [double]"55.21111838815975"
$ws.Range("U9").Value = "[46.59051136950049, 63.831725406819004]"
$ws.Range("V9").Value = [double]"2.220446049250313e-16"
$ws.Range("W9").Value = [double]"2.220446049250313e-16"
$ws.Range("X9").Value = [double]"17.19893893893915"
$ws.Range("Y9").Value = [double]"16.09705705705726"
$ws.Range("Z9").Value = [double]"18.30082082082105"
$ws.Range("F10").Value = [double]"23.88000000000029"
$ws.Range("H10").Value = [double]"3.505958412475252e-09"
$ws.Range("I10").Value = [double]"3.505958412475252e-09"
$ws.Range("L10").Value = [double]"50.59433126266978"
$ws.Range("M10").Value = "[35.65280972761231, 65.53585279772726]"
$ws.Range("N10").Value = [double]"1.884595679513268e-08"
$ws.Range("O10").Value = [double]"1.884595679513268e-08"
$ws.Range("P10").Value = [double]"1.767342413731195"
$ws.Range("Q10").Value = "[1.4277107755052718, 2.1069740519571187]"
$ws.Range("R10").Value = [double]"1.170175067954915e-13"
$ws.Range("S10").Value = [double]"1.170175067954915e-13"
$ws.Range("T10").Value = [double]"50.72481875824636"
$ws.Range("U10").Value = "[41.597442995151034, 59.85219452134168]"
$ws.Range("V10").Value = [double]"1.354472090042691e-14"
$ws.Range("W10").Value = [double]"1.354472090042691e-14"
$ws.Range("X10").Value = [double]"17.16300300300321"
$ws.Range("Y10").Value = [double]"15.87219219219238"
$ws.Range("Z10").Value = [double]"18.45381381381404"
$ws.Range("F11").Value = [double]"23.88000000000029"
$ws.Range("H11").Value = [double]"1.398881011027697e-14"
$ws.Range("I11").Value = [double]"1.398881011027697e-14"
$ws.Range("L11").Value = [double]"61.87204591682443"
$ws.Range("M11").Value = "[49.560982243219655, 74.18310959042921]"
$ws.Range("N11").Value = [double]"3.548272786702e-13"
$ws.Range("O11").Value = [double]"3.548272786702e-13"
$ws.Range("P11").Value = [double]"1.729605565039425"
$ws.Range("Q11").Value = "[1.515763422452732, 1.9434477076261185]"
$ws.Range("T11").Value = [double]"53.67780005594252"
$ws.Range("U11").Value = "[46.27060920979768, 61.084990902087355]"
$ws.Range("X11").Value = [double]"17.30642642642664"
$ws.Range("Y11").Value = [double]"16.4936936936939"
$ws.Range("Z11").Value = [double]"18.11915915915938"
$ws.Range("F12").Value = [double]"23.88000000000029"
$ws.Range("H12").Value = [double]"3.573283668956151e-09"
$ws.Range("I12").Value = [double]"3.573283668956151e-09"
$ws.Range("L12").Value = [double]"51.2115436258944"
$ws.Range("M12").Value = "[35.925224335864826, 66.49786291592397]"
$ws.Range("N12").Value = [double]"2.415200350824875e-08"
$ws.Range("O12").Value = [double]"2.415200350824875e-08"
$ws.Range("P12").Value = [double]"1.729605565039425"
$ws.Range("Q12").Value = "[1.4025528763774249, 2.0566582537014257]"
$ws.Range("R12").Value = [double]"6.94999613415348e-14"
$ws.Range("S12").Value = [double]"6.94999613415348e-14"
$ws.Range("T12").Value = [double]"52.42316628760506"
$ws.Range("U12").Value = "[43.17983206503502, 61.66650051017509]"
$ws.Range("V12").Value = [double]"6.883382752675971e-15"
$ws.Range("W12").Value = [double]"6.883382752675971e-15"
$ws.Range("X12").Value = [double]"17.30642642642664"
$ws.Range("Y12").Value = [double]"16.06342342342362"
$ws.Range("Z12").Value = [double]"18.54942942942966"
$ws.Range("F13").Value = [double]"23.88000000000029"
$ws.Range("H13").Value = [double]"2.098321516541546e-14"
$ws.Range("I13").Value = [double]"2.098321516541546e-14"
$ws.Range("L13").Value = [double]"59.87726817713192"
$ws.Range("M13").Value = "[46.62449844363461, 73.13003791062924]"
$ws.Range("N13").Value = [double]"9.161560399206792e-12"
$ws.Range("O13").Value = [double]"9.161560399206792e-12"
$ws.Range("P13").Value = [double]"1.968605606753964"
$ws.Range("Q13").Value = "[1.7296055650394244, 2.207605648468503]"
$ws.Range("R13").Value = [double]"0"
$ws.Range("S13").Value = [double]"0"
$ws.Range("T13").Value = [double]"53.60351808156867"
$ws.Range("U13").Value = "[46.35110975779574, 60.855926405341606]"
$ws.Range("X13").Value = [double]"16.39807807807829"
$ws.Range("Y13").Value = [double]"15.48972972972992"
$ws.Range("Z13").Value = [double]"17.30642642642665"
$ws.Range("F14").Value = [double]"23.88000000000029"
$ws.Range("H14").Value = [double]"8.70559180299324e-12"
$ws.Range("I14").Value = [double]"8.70559180299324e-12"
$ws.Range("L14").Value = [double]"54.81218505314664"
$ws.Range("M14").Value = "[40.773537896397556, 68.85083220989571]"
$ws.Range("N14").Value = [double]"5.44452705142362e-10"
$ws.Range("O14").Value = [double]"5.44452705142362e-10"
$ws.Range("P14").Value = [double]"1.905710858934349"
$ws.Range("Q14").Value = "[1.6289739685280402, 2.182447749340657]"
$ws.Range("T14").Value = [double]"50.09761166214037"
$ws.Range("U14").Value = "[42.08776597351533, 58.107457350765415]"
$ws.Range("V14").Value = [double]"2.220446049250313e-16"
$ws.Range("W14").Value = [double]"2.220446049250313e-16"
$ws.Range("X14").Value = [double]"16.63711711711732"
$ws.Range("Y14").Value = [double]"15.58534534534554"
$ws.Range("Z14").Value = [double]"17.68888888888911"
